$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the year values in column B from 1975-1994 to 2000-2019
for ($i = 0; $i -lt 20; $i++) {
    $row = 5 + $i
    $ws.Cells.Item($row, 2).Value = 2000 + $i
}

# Update the selection to J11
$ws.Range("J11").Select()
